$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 64 (15.03.2022) : fill in the missing "Fin" time, Tâche ---
$ws.Range("C64").Value2 = 0.70486111111111116
$ws.Range("E64").Value2 = "Refactor"

# --- Row 65 (16.03.2022) : new entry ---
$ws.Range("A65").Value2 = 44636
$ws.Range("B65").Value2 = 0.33333333333333331
$ws.Range("C65").Value2 = 0.34861111111111115
$ws.Range("E65").Value2 = "refactor Astar pathfinding"

# --- Row 66 (16.03.2022) : new entry ---
$ws.Range("A66").Value2 = 44636
$ws.Range("B66").Value2 = 0.34930555555555554
$ws.Range("C66").Value2 = 0.35625000000000001
$ws.Range("E66").Value2 = "implémentation de Astar`nPathfinding"
$ws.Range("F66").Value2 = "PROBLEME : la grille doit s'adapter a une grille plus grande.."

# --- Row 67 (16.03.2022) : new entry ---
$ws.Range("A67").Value2 = 44636
$ws.Range("B67").Value2 = 0.35625000000000001
$ws.Range("C67").Value2 = 0.39930555555555558
$ws.Range("E67").Value2 = "Analyse du problème d'adaptation de grille"
$ws.Range("F67").Value2 = "Problematique : `nSoit on calcule dynamiquement a chaque besoin les cellule concernée`non crée une grille temporaire adapté qui calculera une fois lors du callback tout les cas"

# --- Move the view / active selection down to the newly filled rows ---
$ws.Range("A68").Select()
